$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18
$ws.Cells.Item(18, 1).Value = "WGG 02"
$ws.Cells.Item(18, 2).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(18, 3).Value = "20-01-2026"
$ws.Cells.Item(18, 4).Value = 286962
$ws.Cells.Item(18, 5).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(18, 6).Value = 34400000000
$ws.Cells.Item(18, 7).Value = "NEFT"
$ws.Cells.Item(18, 8).Value = "SBIN0003229"
$ws.Cells.Item(18, 9).Value = "AAAFW8862C"
$ws.Cells.Item(18, 10).Value = "32AAAFW8862C1Z9"
$ws.Cells.Item(18, 11).Value = ""
$ws.Cells.Item(18, 12).Value = "d7253f4f-1311-4793-8650-8fbd14187c4d"
$ws.Cells.Item(18, 13).Value = ""
$ws.Cells.Item(18, 14).Value = ""
$ws.Cells.Item(18, 15).Value = ""
$ws.Cells.Item(18, 16).Value = ""
$ws.Cells.Item(18, 17).Value = ""
$ws.Cells.Item(18, 18).Value = ""
$ws.Cells.Item(18, 19).Value = ""
$ws.Cells.Item(18, 20).Value = ""
$ws.Cells.Item(18, 21).Value = "pending"
$ws.Cells.Item(18, 22).Value = 500
$ws.Cells.Item(18, 23).Value = ""
$ws.Cells.Item(18, 24).Value = "PAYMENT TESTING RPA_UNIQUE_ID : a00c6ae7-4289-4b6d-a5bd-b709ec0e71cc"
$ws.Cells.Item(18, 25).Value = "HO"
$ws.Cells.Item(18, 26).Value = 0
$ws.Cells.Item(18, 27).Value = "midhuncraju12@gmail.com"
$ws.Cells.Item(18, 28).Value = "ESTIMATION NOT MATCHED"
$ws.Cells.Item(18, 29).Value = 0
$ws.Cells.Item(18, 30).Value = 0
$ws.Cells.Item(18, 31).Value = 0
$ws.Cells.Item(18, 32).Value = ""
$ws.Cells.Item(18, 33).Value = ""
$ws.Cells.Item(18, 34).Value = ""
$ws.Cells.Item(18, 35).Value = ""
$ws.Cells.Item(18, 36).Value = ""
$ws.Cells.Item(18, 37).Value = ""
$ws.Cells.Item(18, 38).Value = ""
$ws.Cells.Item(18, 39).Value = ""
$ws.Cells.Item(18, 40).Value = ""
$ws.Cells.Item(18, 41).Value = ""

# Row 19
$ws.Cells.Item(19, 1).Value = "WGG 02"
$ws.Cells.Item(19, 2).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(19, 3).Value = "20-01-2026"
$ws.Cells.Item(19, 4).Value = 286962
$ws.Cells.Item(19, 5).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(19, 6).Value = 34400000000
$ws.Cells.Item(19, 7).Value = "NEFT"
$ws.Cells.Item(19, 8).Value = "SBIN0003229"
$ws.Cells.Item(19, 9).Value = "AAAFW8862C"
$ws.Cells.Item(19, 10).Value = "32AAAFW8862C1Z9"
$ws.Cells.Item(19, 11).Value = ""
$ws.Cells.Item(19, 12).Value = "7794d18d-54dd-47ac-ab4c-8a7e10988366"
$ws.Cells.Item(19, 13).Value = ""
$ws.Cells.Item(19, 14).Value = ""
$ws.Cells.Item(19, 15).Value = ""
$ws.Cells.Item(19, 16).Value = ""
$ws.Cells.Item(19, 17).Value = ""
$ws.Cells.Item(19, 18).Value = ""
$ws.Cells.Item(19, 19).Value = ""
$ws.Cells.Item(19, 20).Value = ""
$ws.Cells.Item(19, 21).Value = "pending"
$ws.Cells.Item(19, 22).Value = 1500
$ws.Cells.Item(19, 23).Value = ""
$ws.Cells.Item(19, 24).Value = "PAYMENT TESTING RPA_UNIQUE_ID : 2a4421c6-2e5e-4d56-9502-2a6f01ddcde5"
$ws.Cells.Item(19, 25).Value = "HO"
$ws.Cells.Item(19, 26).Value = 0
$ws.Cells.Item(19, 27).Value = "midhuncraju12@gmail.com"
$ws.Cells.Item(19, 28).Value = "ESTIMATION NOT MATCHED"
$ws.Cells.Item(19, 29).Value = 0
$ws.Cells.Item(19, 30).Value = 0
$ws.Cells.Item(19, 31).Value = 0
$ws.Cells.Item(19, 32).Value = ""
$ws.Cells.Item(19, 33).Value = ""
$ws.Cells.Item(19, 34).Value = ""
$ws.Cells.Item(19, 35).Value = ""
$ws.Cells.Item(19, 36).Value = ""
$ws.Cells.Item(19, 37).Value = ""
$ws.Cells.Item(19, 38).Value = ""
$ws.Cells.Item(19, 39).Value = ""
$ws.Cells.Item(19, 40).Value = ""
$ws.Cells.Item(19, 41).Value = ""
